$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New session row 13
$ws.Range("A13").Value = 20240523
$ws.Range("B13").Value = 1
$ws.Range("C13").Value = 3
$ws.Range("D13").Value = 2
$ws.Range("E13").Value = 5
$ws.Range("F13").Value = 4
$ws.Range("G13").Value = 6

# Row 14 only has a value in column F
$ws.Range("F14").Value = 27

# New session row 15
$ws.Range("A15").Value = 20240610
$ws.Range("B15").Value = 4
$ws.Range("C15").Value = 6
$ws.Range("D15").Value = 2
$ws.Range("E15").Value = 5
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 3

# Update the active selection to reflect where the user left off entering data
$ws.Range("A16").Select()
